# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for sheet "展览"
$exhibitionUpdates = @{
    2  = 1200
    3  = 425
    5  = 149
    7  = 12364
    11 = 164
    12 = 12194
    13 = 4839
    14 = 4718
    15 = 136
    19 = 953
}

# Row -> new F-column value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 1200
    3  = 425
    5  = 149
    9  = 12364
    13 = 164
    14 = 12194
    15 = 4839
    16 = 4718
    17 = 136
    21 = 953
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
